$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 0: Split column A/B width definition (col A keeps 30.7109375, now solely min=1 max=1)
$ws.Columns.Item(2).ColumnWidth = 59.8

# Step 1: Update Objetivos answer text (row 10) - content only, style/height unchanged
$ws.Range("B10").Value = "O objetivo desta disciplina é apresentar aos alunos fundamentos gerais relativos à obtenção de metais e ligas via processos pirometalúrgicos. Serão feitas aplicações com relação à produção de ferro gusa / aços e cobre por esta via e eventualmente outros metais mais relevantes através da apresentação de trabalhos pelos alunos."
$ws.Range("C10").Value = "O objetivo desta disciplina é apresentar aos alunos fundamentos gerais relativos à obtenção de metais e ligas via processos pirometalúrgicos. Serão feitas aplicações com relação à produção de ferro gusa / aços e cobre por esta via e eventualmente outros metais mais relevantes através da apresentação de trabalhos pelos alunos."

# Step 2: Clear old content & formatting for rows 13-23 (to be rebuilt as rows 13-25)
$ws.Range("A13:C23").Clear()

# Step 3: Populate values for rows 13-25
$ws.Range("B13").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C13").Value = "3577649 - Carlos Angelo Nunes"

$ws.Range("B14").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C14").Value = "1922320 - Sebastiao Ribeiro"

$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "1- Fundamentos de Pirometalurgia; 2- Balanço de massa e energia; 3- Equilíbrio de reações em sistemas envolvendo gases e fases condensadas; 4 - Preparação de matérias-primas; 5- Redução de óxidos, obtenção de gusa em alto-forno; 6- Refino do aço por processos pneumáticos; 7- Escórias; 8- Produção pirometalúrgica do cobre."
$ws.Range("C15").Value = "1- Fundamentos de Pirometalurgia; 2- Balanço de massa e energia; 3- Equilíbrio de reações em sistemas envolvendo gases e fases condensadas; 4 - Preparação de matérias-primas; 5- Redução de óxidos, obtenção de gusa em alto-forno; 6- Refino do aço por processos pneumáticos; 7- Escórias; 8- Produção pirometalúrgica do cobre."

$ws.Range("A16").Value = "Short syllabus:"

$ws.Range("A17").Value = "Programa:"
$ws.Range("B17").Value = "1- Fundamentos de Pirometalurgia: minérios, fluxogramas, características da extração de metais, unidades de medida, estequiometria; 2- Balanço de massa e energia: balanço de materiais, variação de entalpia em reações químicas, balanço de entalpia; 3 - Equilíbrio de reações em sistemas envolvendo gases e fases condensadas: variação de energia livre de Gibbs em reações químicas, constante de equilíbrio, diagrama de Ellingham; 4- Preparação de matérias primas: secagem; calcinação; coqueificação, ustulação de sulfetos, aglomeração de minérios; 5- Redução de óxidos: termodinâmica da redução de óxidos, operação do alto forno, comportamento de impurezas durante a fabricação do gusa; 6- Refino do aço por processos pneumáticos: descarburação, desfosforação, dessulfuração, desoxidação; 7 - Escórias: escória de alto-forno; escória de aciaria LD; 8- Produção pirometalúrgica do cobre: fundamentos, minério, forno de fusão para matte, escória, conversão da matte, forno de conversão, processo de refino."
$ws.Range("C17").Value = "1- Fundamentos de Pirometalurgia: minérios, fluxogramas, características da extração de metais, unidades de medida, estequiometria; 2- Balanço de massa e energia: balanço de materiais, variação de entalpia em reações químicas, balanço de entalpia; 3 - Equilíbrio de reações em sistemas envolvendo gases e fases condensadas: variação de energia livre de Gibbs em reações químicas, constante de equilíbrio, diagrama de Ellingham; 4- Preparação de matérias primas: secagem; calcinação; coqueificação, ustulação de sulfetos, aglomeração de minérios; 5- Redução de óxidos: termodinâmica da redução de óxidos, operação do alto forno, comportamento de impurezas durante a fabricação do gusa; 6- Refino do aço por processos pneumáticos: descarburação, desfosforação, dessulfuração, desoxidação; 7 - Escórias: escória de alto-forno; escória de aciaria LD; 8- Produção pirometalúrgica do cobre: fundamentos, minério, forno de fusão para matte, escória, conversão da matte, forno de conversão, processo de refino."

$ws.Range("A18").Value = "Syllabus:"

$ws.Range("A19").Value = "Avaliação:"

$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "O curso será ministrado na forma de aulas expositivas."
$ws.Range("C20").Value = "O curso será ministrado na forma de aulas expositivas."

$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF=(P1+P2)/2. Poderão ser solicitados trabalhos aos alunos e que irão também compor parte da nota P2."
$ws.Range("C21").Value = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF=(P1+P2)/2. Poderão ser solicitados trabalhos aos alunos e que irão também compor parte da nota P2."

$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Para os alunos que obtiverem 3,0"
$ws.Range("C22").Value = "Para os alunos que obtiverem 3,0"

$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "1) Terkel Rosenqvist. Principles of extractive metallurgy, 2nd edition, McGraw-Hill International Editions - Materials Science & Engineering Series, London, 1983. 2) Carlos Antônio da Silva, Danton Heleno Gameiro e Itavahn Alves da Silva. Balanço de energia em processos metalúrgicos, Escola de Minas - Departamento de Metalurgia, Universidade Federal de Ouro Preto (apostila).3) Fathi Habashi. Extractive Metallurgy, Gordon and Breach Science Publishers, 1986. 4) Alan H. Cottrell. Introdução à metalurgia, 2a edição, Fundação Calouste Gulbenkian, Lisboa, 1975."
$ws.Range("C23").Value = "1) Terkel Rosenqvist. Principles of extractive metallurgy, 2nd edition, McGraw-Hill International Editions - Materials Science & Engineering Series, London, 1983. 2) Carlos Antônio da Silva, Danton Heleno Gameiro e Itavahn Alves da Silva. Balanço de energia em processos metalúrgicos, Escola de Minas - Departamento de Metalurgia, Universidade Federal de Ouro Preto (apostila).3) Fathi Habashi. Extractive Metallurgy, Gordon and Breach Science Publishers, 1986. 4) Alan H. Cottrell. Introdução à metalurgia, 2a edição, Fundação Calouste Gulbenkian, Lisboa, 1975."

$ws.Range("A24").Value = "Requisitos:"

$ws.Range("B25").Value = "LOM3015 -  Termodinâmica de Materiais  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOM3015 -  Termodinâmica de Materiais  (Requisito fraco)`n"

# Step 4: Apply cell styles by copying formats from template rows (A=bold label, B=normal wrap, C=red wrap)
$excel.CutCopyMode = $false
$ws.Range("B3:C3").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B3:C3").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A3:A3").Copy()
$ws.Range("A16:A16").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A3:A3").Copy()
$ws.Range("A18:A18").PasteSpecial(-4122)
$ws.Range("A3:A3").Copy()
$ws.Range("A19:A19").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A23:C23").PasteSpecial(-4122)
$ws.Range("A3:A3").Copy()
$ws.Range("A24:A24").PasteSpecial(-4122)
$ws.Range("B3:C3").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 5: Set row heights exactly, and reset rows without custom height to default via AutoFit
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(25).RowHeight = 30
